$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# The "Security sensors monitoring" and "Contact police" backlog items (originally
# rows 7 and 8) are re-prioritised to the top of the backlog, becoming rows 4 and 5.
# The three rows that used to sit at 4-6 (sprinkler settings, water usage report,
# electric usage report) slide down to rows 6-8. Column A (the running priority
# number) and the row heights are unaffected - only the User Story (B) and Story
# Points (C) content moves.
$rows = @(4, 5, 6, 7, 8)
$newStoryText = @(
    "As an User I want the SHAS system to monitor security sensors so that the system reports any unauthorized access",
    "As an User I want the SHAS system to contact the police if security sensors a tripped for my safety",
    "As an User I want ot be able to set sprinkler settings from the SHAS system to reduce water usage",
    "As an User I want to view water usage reports to be able to monitor water usage ",
    "As an User I want to view Electric usage report to be able to monitor electricity usage"
)
$newPoints = @(10, 10, 2, 5, 5)
# Row 6 (the sprinkler-settings story) keeps the distinct "left aligned" look that
# used to belong to row 4; every other row in the block uses the plain wrap style.
$leftAlignedRows = @(6)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $bCell = $ws1.Cells.Item($r, 2)
    $bCell.Value2 = $newStoryText[$i]
    $cCell = $ws1.Cells.Item($r, 3)
    $cCell.Value2 = $newPoints[$i]

    if ($leftAlignedRows -contains $r) {
        $bCell.HorizontalAlignment = -4131
    } else {
        $bCell.HorizontalAlignment = 1
    }
}

# Add the new (empty) Sheet2 right after Sheet1 - this is the "Sprint Review Power
# Point" placeholder tab mentioned in the commit message.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Match the workbook's usual (1"/0.75"/0.5") page-setup instead of the host's
# narrower defaults.
$ps2 = $ws2.PageSetup
$ps2.LeftMargin = 54
$ps2.RightMargin = 54
$ps2.TopMargin = 72
$ps2.BottomMargin = 72
$ps2.HeaderMargin = 36
$ps2.FooterMargin = 36

# Restore the active sheet/selection as recorded in the saved workbook.
$ws1.Activate()
$ws1.Range("A1:C6").Select()
